$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.090940713882446
$ws.Range("B1").Value = 2.397950172424316
$ws.Range("C1").Value = 3.007592916488647
$ws.Range("D1").Value = 5.951735973358154
$ws.Range("E1").Value = 2.499330282211304
